$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("LEFT RIGHT MID Functions")
$dst = $wb.Worksheets.Item("SEARCH Function")
$srcCell = $src.Cells.Item(3, 5)   # E3 has s="17": numFmtId 0 fontId3 borderId3
$srcCell.Copy()
$dstRange = $dst.Range("A2:A6")
$dstRange.PasteSpecial(-4122)   # xlPasteFormats
